$p = $ppt.ActivePresentation

# --- Slide 4: "Intravenous Catheter in Peripheral Vein ("IV")" ---
# Content Placeholder 2 (Shape 2): replace last bullet text.
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$tr4.Text = "IV catheter placed into a vein in the hand or arm`rAllows administration of chemotherapy and fluids`rPlaced at the beginning of each dose`rRemoved that day at the end of treatment`rNot suitable for FLOT chemotherapy"

# --- Slide 6: "Central Venous Port" ---
# Content Placeholder 2 (Shape 2): insert "May shower within 24 hours" bullet
# before "No special care at home", and drop the "Allows showering, bathing,
# swimming" bullet.
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange
$tr6.Text = "Implantable device that makes the administration of chemotherapy easier`rMay shower within 24 hours`rNo special care at home`rSuitable for FLOT chemotherapy`rAllows for blood draws"
